$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-10-28 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-10-29 Sunday", 2)

# Update the table of division problems (Table 1), row by row, cell by cell.
$tbl = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="88÷6=14, 4"},
    @{Row=1;  Col=2; Text="24÷9=2, 6"},
    @{Row=1;  Col=3; Text="93÷6=15, 3"},
    @{Row=1;  Col=4; Text="59÷6=9, 5"},
    @{Row=1;  Col=5; Text="49÷6=8, 1"},

    @{Row=5;  Col=1; Text="99÷7=14, 1"},
    @{Row=5;  Col=2; Text="42÷9=4, 6"},
    @{Row=5;  Col=3; Text="51÷2=25, 1"},
    @{Row=5;  Col=4; Text="86÷5=17, 1"},
    @{Row=5;  Col=5; Text="38÷2=19, 0"},

    @{Row=9;  Col=1; Text="62÷6=10, 2"},
    @{Row=9;  Col=2; Text="96÷8=12, 0"},
    @{Row=9;  Col=3; Text="60÷9=6, 6"},
    @{Row=9;  Col=4; Text="30÷8=3, 6"},
    @{Row=9;  Col=5; Text="80÷9=8, 8"},

    @{Row=13; Col=1; Text="81÷4=20, 1"},
    @{Row=13; Col=2; Text="80÷6=13, 2"},
    @{Row=13; Col=3; Text="61÷6=10, 1"},
    @{Row=13; Col=4; Text="59÷6=9, 5"},
    @{Row=13; Col=5; Text="46÷4=11, 2"},

    @{Row=17; Col=1; Text="54÷5=10, 4"},
    @{Row=17; Col=2; Text="47÷9=5, 2"},
    @{Row=17; Col=3; Text="30÷4=7, 2"},
    @{Row=17; Col=4; Text="53÷5=10, 3"},
    @{Row=17; Col=5; Text="50÷6=8, 2"}
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $u.Text
}
